# Generate Report for Handoff
# Swaps the "74628366..." and "a43aa681..." file rows so that the
# 74628366 file moves from "In Translation" to "Ready for handoff"
# (new priority "mt" and new handoff timestamps), while the a43aa681
# file keeps its previous "In Translation" data.

$wb = $excel.ActiveWorkbook

$urlFor74628366 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/df1b57ddfd6f508de9494e42f1137b57815904be/e2e/74628366-9c09-4cb5-b3c2-e3b52ee9ed79.md"
$urlForA43aa681 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/df1b57ddfd6f508de9494e42f1137b57815904be/e2e/a43aa681-4ac8-4762-98c1-586d9318e933.md"

# ---------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------
$ws = $wb.Worksheets("Overview")

$ws.Range("A2").Value = "a43aa681-4ac8-4762-98c1-586d9318e933.md"
$ws.Range("B2").Value = "e2e\a43aa681-4ac8-4762-98c1-586d9318e933.md"
$ws.Range("E2").Value = "In Translation"
$ws.Range("F2").Value = "In Translation"
$ws.Range("G2").Value = "2016-08-27 12:12:24"

$ws.Range("A3").Value = "74628366-9c09-4cb5-b3c2-e3b52ee9ed79.md"
$ws.Range("B3").Value = "e2e\74628366-9c09-4cb5-b3c2-e3b52ee9ed79.md"
$ws.Range("E3").Value = "Ready for handoff"
$ws.Range("F3").Value = "Ready for handoff"
$ws.Range("G3").Value = "2016-08-27 12:12:55"

# Recreate the hyperlinks, keeping relationship order/targets stable
# (rId2 -> 74628366 url, rId3 -> a43aa681 url) but repointing the
# display text to the new cell content.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), $urlFor74628366, "", "", "e2e\a43aa681-4ac8-4762-98c1-586d9318e933.md")
$ws.Hyperlinks.Add($ws.Range("B3"), $urlForA43aa681, "", "", "e2e\74628366-9c09-4cb5-b3c2-e3b52ee9ed79.md")

# Widen the status columns (E, F) to fit "Ready for handoff"
$ws.Columns.Item(5).ColumnWidth = 16.333333333333332
$ws.Columns.Item(6).ColumnWidth = 16.333333333333332

# ---------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------
$ws = $wb.Worksheets("zh-cn")

$ws.Range("A2").Value = "a43aa681-4ac8-4762-98c1-586d9318e933.md"
$ws.Range("C2").Value = "In Translation"
$ws.Range("E2").Value = "ht"
$ws.Range("G2").Value = "a43aa681-4ac8-4762-98c1-586d9318e933.2756ec0c465a47e6967c66db16d3e1ec21402547.zh-cn.xlf"
$ws.Range("H2").Value = "2016-08-27 12:12:20"

$ws.Range("A3").Value = "74628366-9c09-4cb5-b3c2-e3b52ee9ed79.md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("E3").Value = "mt"
$ws.Range("G3").Value = "74628366-9c09-4cb5-b3c2-e3b52ee9ed79.e9210b563fbf3a2432cf8a03b28a77ed6f38669f.zh-cn.xlf"
$ws.Range("H3").Value = "2016-08-27 12:12:51"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $urlFor74628366, "", "", "a43aa681-4ac8-4762-98c1-586d9318e933.md")
$ws.Hyperlinks.Add($ws.Range("A3"), $urlForA43aa681, "", "", "74628366-9c09-4cb5-b3c2-e3b52ee9ed79.md")

$ws.Columns.Item(3).ColumnWidth = 16.333333333333332

# ---------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------
$ws = $wb.Worksheets("de-de")

$ws.Range("A2").Value = "a43aa681-4ac8-4762-98c1-586d9318e933.md"
$ws.Range("C2").Value = "In Translation"
$ws.Range("E2").Value = "ht"
$ws.Range("G2").Value = "a43aa681-4ac8-4762-98c1-586d9318e933.2756ec0c465a47e6967c66db16d3e1ec21402547.de-de.xlf"
$ws.Range("H2").Value = "2016-08-27 12:12:24"

$ws.Range("A3").Value = "74628366-9c09-4cb5-b3c2-e3b52ee9ed79.md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("E3").Value = "mt"
$ws.Range("G3").Value = "74628366-9c09-4cb5-b3c2-e3b52ee9ed79.e9210b563fbf3a2432cf8a03b28a77ed6f38669f.de-de.xlf"
$ws.Range("H3").Value = "2016-08-27 12:12:55"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $urlFor74628366, "", "", "a43aa681-4ac8-4762-98c1-586d9318e933.md")
$ws.Hyperlinks.Add($ws.Range("A3"), $urlForA43aa681, "", "", "74628366-9c09-4cb5-b3c2-e3b52ee9ed79.md")

$ws.Columns.Item(3).ColumnWidth = 16.333333333333332

Write-Host "Report updated for handoff."
